$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("債務")

# --- Fix header row (row 1) ---
# The header row previously (erroneously) duplicated the first data row's
# values.  Replace it with proper column headers, matching the pattern used
# on the other sheets (species/debtor/owner/total/... + the common
# legislator/source metadata columns).
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Give the new header cells (H1:N1) the same bold/bordered style as the
# existing header cells (B1:G1).
$ws.Range("B1:G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# --- Add the missing metadata columns to each data row ---
$rows = @(
    @{ Row = 2; Index = 104 },
    @{ Row = 3; Index = 105 },
    @{ Row = 4; Index = 106 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 8).Value  = "debt"            # H: property_category
    $ws.Cells.Item($row, 9).Value  = "normal"           # I: category

    # Force the date column to be stored as plain text (matching the other
    # sheets) instead of letting Excel auto-convert it to a date serial.
    $ws.Cells.Item($row, 10).NumberFormat = "@"
    $ws.Cells.Item($row, 10).Value = "2012-04-19"       # J: date

    $ws.Cells.Item($row, 11).Value = "李慶華"            # K: legislator_name
    $ws.Cells.Item($row, 12).Value = 607                # L: legislator_id
    $ws.Cells.Item($row, 13).Value = "tmpe2cb1"         # M: source_file
    $ws.Cells.Item($row, 14).Value = $r.Index            # N: index

    $ws.Range("B$row`:G$row").Copy()
    $ws.Range("H$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
